# "using role name ot unit code to assign users to roles"
#
# Every row in column A (Role_Name) now points at the actual role name
# ("مدير اجازه دراسيه بمرتب") instead of the unit code ("aa" /
# "اعاره من الوزاره") that used to be mixed in. Column B (User) keeps the
# same users, they just all now line up under the single role name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$roleName = "مدير اجازه دراسيه بمرتب"

$ws.Range("A1").Value = $roleName
$ws.Range("A2").Value = $roleName
$ws.Range("A3").Value = $roleName
$ws.Range("A4").Value = $roleName
$ws.Range("A5").Value = $roleName
$ws.Range("A6").Value = $roleName

# Column B (users) is untouched - user3/mfa_remon/mfa_omar/youssef1/mf1/mfa_salam
# stay exactly where they were.

# Selection moved from B6 to A6 as part of the edit.
$ws.Range("A6").Select()

# Page orientation explicitly set to portrait.
$ws.PageSetup.Orientation = 1
